# Add a reviewer comment on the "True" answer to the second true/false
# question ("Event driven programming is a more advanced paradigm than
# Object-Oriented programming"), correcting it to False.

$word.UserName = "Amir Gamil"
$word.UserInitials = "AG"

$d = $word.ActiveDocument

# Locate the paragraph that consists solely of the answer "True" - it is
# the response directly under the second True/False question. (Cast to
# [string] explicitly since a bare paragraph of "True"/"False" text can
# otherwise be auto-coerced to a PowerShell boolean.)
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = [string]$para.Range.Text
    $trimmed = $text.TrimEnd([char]13, [char]7)
    if ($trimmed -eq "True") {
        $target = $para
        break
    }
}

$r = $target.Range
# Drop the trailing paragraph mark so the comment anchors only the word.
[void]$r.MoveEnd(1, -1)

$comment = $d.Comments.Add($r, "False.  It is just another paradigm and can co-exist with object oriented programming.")
